$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. add_new_document_type: edit existing row + append two new document-type
#    rows (Advanced level certiciate / Diploma certiciate).
# ---------------------------------------------------------------------------
$addNew = $wb.Worksheets.Item("add_new_document_type")
$addNew.Range("A2").Value2 = "Auto_DT_004"
$addNew.Range("A3").Value2 = "Auto_DT_005"
$addNew.Range("B3").Value2 = "Advanced level certiciate"
$addNew.Range("C3").Value2 = "Active"
$addNew.Range("D3").Value2 = "Y"
$addNew.Range("A4").Value2 = "Auto_DT_006"
$addNew.Range("B4").Value2 = "Diploma certiciate"
$addNew.Range("C4").Value2 = "Inactive"
$addNew.Range("D4").Value2 = "Y"
$addNew.Range("A2:A4").Select()

# ---------------------------------------------------------------------------
# 2. New sheet "filter_document_type": clone the existing filter_* sheet
#    pattern (filter_awarding_institute) and drop it right after
#    edit_document_type / before search_departments.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("filter_awarding_institute")
$template.Copy($null, $wb.Worksheets.Item("edit_document_type"))
$newSheet = $wb.Worksheets.Item("filter_awarding_institute (2)")
$newSheet.Name = "filter_document_type"
$newSheet.Range("B2").Value2 = "Auto_DT_"
$newSheet.Range("B3").Value2 = "Dummy"

# Match the greyed-out "Dummy" placeholder styling used elsewhere in the
# workbook (search_locations!A3) without introducing a brand-new style.
$wb.Worksheets.Item("search_locations").Range("A3").Copy()
$newSheet.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("I15").Select()

# ---------------------------------------------------------------------------
# 3. filter_awarding_institute: selection now spans the whole filter block
#    with no distinct active cell.
# ---------------------------------------------------------------------------
$template.Range("A1:C4").Select()

# ---------------------------------------------------------------------------
# 4. edit_document_type becomes the active / selected tab.
# ---------------------------------------------------------------------------
$editDocType = $wb.Worksheets.Item("edit_document_type")
$editDocType.Activate()
